# Adding Ammonia, Methanol, Jet Fuel into the Trade Links
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

# Data layout (columns B..H): Reg1, Reg2, Comm, Comm1, Comm2, Tech, TradeLink
# Existing rows 4-11 cover commodity H2 for regions DKISLBH/DKISL1/DKISL2/DKISL3.
# New rows 12-35 repeat the same region pattern for three more commodities:
#   MOE (Methanol) rows 12-19, KRE (Jet Fuel) rows 20-27, AMM (Ammonia) rows 28-35

$regions = @(
    @{ Reg1 = "DKISLBH"; Reg2 = "DKE"; Suffix = "01" },
    @{ Reg1 = "DKISLBH"; Reg2 = "DKE"; Suffix = "02" },
    @{ Reg1 = "DKISL1";  Reg2 = "DKW"; Suffix = "01" },
    @{ Reg1 = "DKISL1";  Reg2 = "DKW"; Suffix = "02" },
    @{ Reg1 = "DKISL2";  Reg2 = "DKW"; Suffix = "01" },
    @{ Reg1 = "DKISL2";  Reg2 = "DKW"; Suffix = "02" },
    @{ Reg1 = "DKISL3";  Reg2 = "DKW"; Suffix = "01" },
    @{ Reg1 = "DKISL3";  Reg2 = "DKW"; Suffix = "02" }
)

# Reg1 / Reg2 columns, base row 12 through 35
$row = 12
for ($block = 0; $block -lt 3; $block++) {
    foreach ($reg in $regions) {
        $ws.Cells.Item($row, 2).Value = $reg.Reg1
        $ws.Cells.Item($row, 3).Value = $reg.Reg2
        $row = $row + 1
    }
}

# Comm / Comm1 / Comm2 columns: MOE block first, then KRE block
$row = 12
foreach ($reg in $regions) {
    $ws.Cells.Item($row, 4).Value = "MOE"
    $ws.Cells.Item($row, 5).Value = "MOE"
    $ws.Cells.Item($row, 6).Value = "MOE"
    $row = $row + 1
}
foreach ($reg in $regions) {
    $ws.Cells.Item($row, 4).Value = "KRE"
    $ws.Cells.Item($row, 5).Value = "KRE"
    $ws.Cells.Item($row, 6).Value = "KRE"
    $row = $row + 1
}

# TradeLink column (G): MOE block, then KRE block
$row = 12
foreach ($reg in $regions) {
    $ws.Cells.Item($row, 7).Value = "TB_MOE_" + $reg.Reg1 + "_" + $reg.Reg2 + "_" + $reg.Suffix
    $row = $row + 1
}
foreach ($reg in $regions) {
    $ws.Cells.Item($row, 7).Value = "TB_KRE_" + $reg.Reg1 + "_" + $reg.Reg2 + "_" + $reg.Suffix
    $row = $row + 1
}

# AMM block: Comm/Comm1/Comm2/TradeLink filled together, row by row
foreach ($reg in $regions) {
    $ws.Cells.Item($row, 4).Value = "AMM"
    $ws.Cells.Item($row, 5).Value = "AMM"
    $ws.Cells.Item($row, 6).Value = "AMM"
    $ws.Cells.Item($row, 7).Value = "TB_AMM_" + $reg.Reg1 + "_" + $reg.Reg2 + "_" + $reg.Suffix
    $row = $row + 1
}

# Tech column (H): "B" for every new row (reuses existing shared string)
for ($r = 12; $r -le 35; $r++) {
    $ws.Cells.Item($r, 8).Value = "B"
}

$ws.Range("G12:G35").Select()
